$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.178962333333333
$ws.Range("H2").Value = 3.536887
$ws.Range("I2").Value = 0.001182125215344215
$ws.Range("J2").Value = 0.001182125215344214
$ws.Range("M2").Value = 0.366183
$ws.Range("N2").Value = 1.098549
$ws.Range("O2").Value = 0.0639836884691917
$ws.Range("P2").Value = 0.0639836884691917
$ws.Range("Q2").Value = 0.431715964107
$ws.Range("R2").Value = 3.885443676963
$ws.Range("S2").Value = 0.00007563673151016038
$ws.Range("T2").Value = 0.00007563673151016035

$ws.Range("G3").Value = 1.178962333333333
$ws.Range("H3").Value = 3.536887
$ws.Range("I3").Value = 0.001182125215344215
$ws.Range("J3").Value = 0.001182125215344214
$ws.Range("O3").Value = 0.2777364052521014
$ws.Range("P3").Value = 0.2777364052521014
$ws.Range("Q3").Value = 1.873965737670111
$ws.Range("R3").Value = 16.865691639031
$ws.Range("S3").Value = 0.0003283192078675685
$ws.Range("T3").Value = 0.0003283192078675684

$ws.Range("G4").Value = 1.178962333333333
$ws.Range("H4").Value = 3.536887
$ws.Range("I4").Value = 0.001182125215344215
$ws.Range("J4").Value = 0.001182125215344214
$ws.Range("O4").Value = 0.6582799062787069
$ws.Range("P4").Value = 0.6582799062787069
$ws.Range("Q4").Value = 4.441599901328223
$ws.Range("R4").Value = 39.974399111954
$ws.Range("S4").Value = 0.0007781692759664859
$ws.Range("T4").Value = 0.0007781692759664856

$ws.Range("I5").Value = 0.9532080272144655
$ws.Range("J5").Value = 0.9532080272144653
$ws.Range("M5").Value = 0.366183
$ws.Range("N5").Value = 1.098549
$ws.Range("O5").Value = 0.0639836884691917
$ws.Range("P5").Value = 0.0639836884691917
$ws.Range("Q5").Value = 348.114664268961
$ws.Range("R5").Value = 3133.031978420649
$ws.Range("S5").Value = 0.06098976545962316
$ws.Range("T5").Value = 0.06098976545962315

$ws.Range("I6").Value = 0.9532080272144655
$ws.Range("J6").Value = 0.9532080272144653
$ws.Range("O6").Value = 0.2777364052521014
$ws.Range("P6").Value = 0.2777364052521014
$ws.Range("S6").Value = 0.2647405709359929
$ws.Range("T6").Value = 0.2647405709359928

$ws.Range("I7").Value = 0.9532080272144655
$ws.Range("J7").Value = 0.9532080272144653
$ws.Range("O7").Value = 0.6582799062787069
$ws.Range("P7").Value = 0.6582799062787069
$ws.Range("S7").Value = 0.6274776908188495
$ws.Range("T7").Value = 0.6274776908188493

$ws.Range("G8").Value = 45.48781433333333
$ws.Range("I8").Value = 0.04560984757019037
$ws.Range("J8").Value = 0.04560984757019036
$ws.Range("M8").Value = 0.366183
$ws.Range("N8").Value = 1.098549
$ws.Range("O8").Value = 0.0639836884691917
$ws.Range("P8").Value = 0.0639836884691917
$ws.Range("Q8").Value = 16.656864316023
$ws.Range("R8").Value = 149.911778844207
$ws.Range("S8").Value = 0.002918286278058381
$ws.Range("T8").Value = 0.00291828627805838

$ws.Range("G9").Value = 45.48781433333333
$ws.Range("I9").Value = 0.04560984757019037
$ws.Range("J9").Value = 0.04560984757019036
$ws.Range("O9").Value = 0.2777364052521014
$ws.Range("P9").Value = 0.2777364052521014
$ws.Range("Q9").Value = 72.30307799669544
$ws.Range("R9").Value = 650.727701970259
$ws.Range("S9").Value = 0.01266751510824097
$ws.Range("T9").Value = 0.01266751510824096

$ws.Range("G10").Value = 45.48781433333333
$ws.Range("I10").Value = 0.04560984757019037
$ws.Range("J10").Value = 0.04560984757019036
$ws.Range("O10").Value = 0.6582799062787069
$ws.Range("P10").Value = 0.6582799062787069
$ws.Range("S10").Value = 0.03002404618389103
$ws.Range("T10").Value = 0.03002404618389102
